$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 14 / column B: task is now done (was "offen", now "done")
# Copy the formatting from B17 (already styled as "done") onto B14, then set its value.
$ws.Range("B17").Copy() | Out-Null
$ws.Range("B14").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$ws.Range("B14").Value = "done"

# New row 18: shopping cart persistence task, marked done
$ws.Range("A18").Value = "Einkaufswagen in DB speichern"
$ws.Range("B17").Copy() | Out-Null
$ws.Range("B18").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$ws.Range("B18").Value = "done"

# New row 19: logging task, still open
$ws.Range("A19").Value = "Logs überall einfügen, allgemeine Fehlerprevention/Abfangen der Fehler"
$ws.Range("B4").Copy() | Out-Null
$ws.Range("B19").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$ws.Range("B19").Value = "offen"

$excel.CutCopyMode = 0

$ws.Range("B19").Select() | Out-Null
